# "Actualiza excel de medidas"
# Add a MAX() summary formula below the measurements table and leave the
# selection on D10 (the last data cell), matching what Excel does when a
# user clicks into D16 and types the formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 16: D16 = MAX(D2:D10)
$ws.Range("D16").Formula = "=MAX(D2:D10)"

# Stamp the (already-default) row height explicitly, as Excel does for a
# freshly-typed-into row.
$ws.Rows.Item(16).RowHeight = 15.75

# Leave the selection on D10.
$ws.Range("D10").Select()
